# Update the model-results sheets (arbolts, bosquets, knnts) with the new
# train_test_split values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("arbolts")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 0.6977341098772637
$ws.Range("C2").Value = 0.8658933832837197
$ws.Range("D2").Value = 0.9305339237683491
$ws.Range("E2").Value = -1.732651912978245

$ws = $wb.Worksheets.Item("bosquets")
$ws.Range("A2").Value = 106
$ws.Range("B2").Value = 0.5045252397713184
$ws.Range("C2").Value = 0.484406475399015
$ws.Range("D2").Value = 0.6959931575806008
$ws.Range("E2").Value = -0.5287266391137646

$ws = $wb.Worksheets.Item("knnts")
$ws.Range("A2").Value = 33
$ws.Range("B2").Value = 0.4302434103874868
$ws.Range("C2").Value = 0.3679452024441739
$ws.Range("D2").Value = 0.6065848682947621
$ws.Range("E2").Value = -0.1611893343233775
